# Deploy the implementation guide.
# - Update the CodeSystem "Status" from active -> draft
# - Update the CodeSystem "Date" to the new publication timestamp
# - Ensure the wrap/alignment formatting (vertical=top, wrap text) used by
#   the header and body styles is (re)applied on every cell that uses it.

$wb = $excel.ActiveWorkbook
$metadata = $wb.Worksheets.Item("Metadata")

# --- Cell value updates -----------------------------------------------
# Status: active -> draft
$metadata.Range("B6").Value = "draft"
# Date: refreshed publication timestamp
$metadata.Range("B8").Value = "2023-08-01T16:12:28+00:00"

# --- Re-apply alignment (vertical top + wrap text) across both sheets -
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $used.VerticalAlignment = -4160
    $used.WrapText = $true
}
